$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.76613666666667
$ws.Range("H2").Value = 53.29841
$ws.Range("I2").Value = 0.7872390387208499
$ws.Range("J2").Value = 0.7872390387208499
$ws.Range("M2").Value = 4.877755666666666
$ws.Range("N2").Value = 14.633267
$ws.Range("O2").Value = 0.09961167132870688
$ws.Range("P2").Value = 0.0996116713287069
$ws.Range("Q2").Value = 86.65887380060778
$ws.Range("R2").Value = 779.92986420547
$ws.Range("S2").Value = 0.07841819638218844
$ws.Range("T2").Value = 0.07841819638218846
$ws.Range("G3").Value = 17.76613666666667
$ws.Range("H3").Value = 53.29841
$ws.Range("I3").Value = 0.7872390387208499
$ws.Range("J3").Value = 0.7872390387208499
$ws.Range("O3").Value = 0.1360673938501395
$ws.Range("P3").Value = 0.1360673938501395
$ws.Range("Q3").Value = 118.3741518915633
$ws.Range("R3").Value = 1065.36736702407
$ws.Range("S3").Value = 0.1071175643358351
$ws.Range("T3").Value = 0.1071175643358351
$ws.Range("G4").Value = 17.76613666666667
$ws.Range("H4").Value = 53.29841
$ws.Range("I4").Value = 0.7872390387208499
$ws.Range("J4").Value = 0.7872390387208499
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 0.246459
$ws.Range("N4").Value = 0.739377
$ws.Range("O4").Value = 0.00503309197542868
$ws.Range("P4").Value = 0.00503309197542868
$ws.Range("Q4").Value = 4.37862427673
$ws.Range("R4").Value = 39.40761849057
$ws.Range("S4").Value = 0.003962246488530097
$ws.Range("T4").Value = 0.003962246488530097
$ws.Range("G5").Value = 17.76613666666667
$ws.Range("H5").Value = 53.29841
$ws.Range("I5").Value = 0.7872390387208499
$ws.Range("J5").Value = 0.7872390387208499
$ws.Range("M5").Value = 37.01331466666667
$ws.Range("N5").Value = 111.039944
$ws.Range("O5").Value = 0.7558718368280999
$ws.Range("P5").Value = 0.7558718368280999
$ws.Range("Q5").Value = 657.5836068543379
$ws.Range("R5").Value = 5918.25246168904
$ws.Range("S5").Value = 0.5950518182207164
$ws.Range("T5").Value = 0.5950518182207164
$ws.Range("G6").Value = 17.76613666666667
$ws.Range("H6").Value = 53.29841
$ws.Range("I6").Value = 0.7872390387208499
$ws.Range("J6").Value = 0.7872390387208499
$ws.Range("M6").Value = 0.167274
$ws.Range("N6").Value = 0.501822
$ws.Range("O6").Value = 0.00341600601762507
$ws.Range("P6").Value = 0.00341600601762507
$ws.Range("Q6").Value = 2.97181274478
$ws.Range("R6").Value = 26.74631470302
$ws.Range("S6").Value = 0.002689213293579799
$ws.Range("T6").Value = 0.002689213293579799
$ws.Range("I7").Value = 0.03648413815195897
$ws.Range("J7").Value = 0.03648413815195897
$ws.Range("M7").Value = 4.877755666666666
$ws.Range("N7").Value = 14.633267
$ws.Range("O7").Value = 0.09961167132870688
$ws.Range("P7").Value = 0.0996116713287069
$ws.Range("Q7").Value = 4.016155409380889
$ws.Range("R7").Value = 36.145398684428
$ws.Range("S7").Value = 0.003634245978304072
$ws.Range("T7").Value = 0.003634245978304073
$ws.Range("I8").Value = 0.03648413815195897
$ws.Range("J8").Value = 0.03648413815195897
$ws.Range("O8").Value = 0.1360673938501395
$ws.Range("P8").Value = 0.1360673938501395
$ws.Range("S8").Value = 0.004964301595205502
$ws.Range("T8").Value = 0.004964301595205503
$ws.Range("I9").Value = 0.03648413815195897
$ws.Range("J9").Value = 0.03648413815195897
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 0.246459
$ws.Range("N9").Value = 0.739377
$ws.Range("O9").Value = 0.00503309197542868
$ws.Range("P9").Value = 0.00503309197542868
$ws.Range("Q9").Value = 0.202924810852
$ws.Range("R9").Value = 1.826323297668
$ws.Range("S9").Value = 0.000183628022963056
$ws.Range("T9").Value = 0.000183628022963056
$ws.Range("I10").Value = 0.03648413815195897
$ws.Range("J10").Value = 0.03648413815195897
$ws.Range("M10").Value = 37.01331466666667
$ws.Range("N10").Value = 111.039944
$ws.Range("O10").Value = 0.7558718368280999
$ws.Range("P10").Value = 0.7558718368280999
$ws.Range("Q10").Value = 30.47533211503289
$ws.Range("R10").Value = 274.2779890352959
$ws.Range("S10").Value = 0.02757733252001139
$ws.Range("T10").Value = 0.02757733252001139
$ws.Range("I11").Value = 0.03648413815195897
$ws.Range("J11").Value = 0.03648413815195897
$ws.Range("M11").Value = 0.167274
$ws.Range("N11").Value = 0.501822
$ws.Range("O11").Value = 0.00341600601762507
$ws.Range("P11").Value = 0.00341600601762507
$ws.Range("Q11").Value = 0.137726943672
$ws.Range("R11").Value = 1.239542493048
$ws.Range("S11").Value = 0.0001246300354749563
$ws.Range("T11").Value = 0.0001246300354749563
$ws.Range("G12").Value = 3.885299333333334
$ws.Range("H12").Value = 11.655898
$ws.Range("I12").Value = 0.1721623203571791
$ws.Range("J12").Value = 0.172162320357179
$ws.Range("M12").Value = 4.877755666666666
$ws.Range("N12").Value = 14.633267
$ws.Range("O12").Value = 0.09961167132870688
$ws.Range("P12").Value = 0.0996116713287069
$ws.Range("Q12").Value = 18.95154083986289
$ws.Range("R12").Value = 170.563867558766
$ws.Range("S12").Value = 0.01714937647060686
$ws.Range("T12").Value = 0.01714937647060686
$ws.Range("G13").Value = 3.885299333333334
$ws.Range("H13").Value = 11.655898
$ws.Range("I13").Value = 0.1721623203571791
$ws.Range("J13").Value = 0.172162320357179
$ws.Range("O13").Value = 0.1360673938501395
$ws.Range("P13").Value = 0.1360673938501395
$ws.Range("Q13").Value = 25.88739589576067
$ws.Range("R13").Value = 232.986563061846
$ws.Range("S13").Value = 0.02342567825019417
$ws.Range("T13").Value = 0.02342567825019417
$ws.Range("G14").Value = 3.885299333333334
$ws.Range("H14").Value = 11.655898
$ws.Range("I14").Value = 0.1721623203571791
$ws.Range("J14").Value = 0.172162320357179
$ws.Range("K14").Value = 3.0
$ws.Range("L14").Value = 1.0
$ws.Range("M14").Value = 0.246459
$ws.Range("N14").Value = 0.739377
$ws.Range("O14").Value = 0.00503309197542868
$ws.Range("P14").Value = 0.00503309197542868
$ws.Range("Q14").Value = 0.957566988394
$ws.Range("R14").Value = 8.618102895546
$ws.Range("S14").Value = 0.0008665087930608996
$ws.Range("T14").Value = 0.0008665087930608995
$ws.Range("G15").Value = 3.885299333333334
$ws.Range("H15").Value = 11.655898
$ws.Range("I15").Value = 0.1721623203571791
$ws.Range("J15").Value = 0.172162320357179
$ws.Range("M15").Value = 37.01331466666667
$ws.Range("N15").Value = 111.039944
$ws.Range("O15").Value = 0.7558718368280999
$ws.Range("P15").Value = 0.7558718368280999
$ws.Range("Q15").Value = 143.8078067988569
$ws.Range("R15").Value = 1294.270261189712
$ws.Range("S15").Value = 0.1301326493209687
$ws.Range("T15").Value = 0.1301326493209687
$ws.Range("G16").Value = 3.885299333333334
$ws.Range("H16").Value = 11.655898
$ws.Range("I16").Value = 0.1721623203571791
$ws.Range("J16").Value = 0.172162320357179
$ws.Range("M16").Value = 0.167274
$ws.Range("N16").Value = 0.501822
$ws.Range("O16").Value = 0.00341600601762507
$ws.Range("P16").Value = 0.00341600601762507
$ws.Range("Q16").Value = 0.6499095606840001
$ws.Range("R16").Value = 5.849186046156
$ws.Range("S16").Value = 0.0005881075223484188
$ws.Range("T16").Value = 0.0005881075223484187
$ws.Range("G17").Value = 0.09285466666666665
$ws.Range("H17").Value = 0.278564
$ws.Range("I17").Value = 0.004114502770011991
$ws.Range("J17").Value = 0.004114502770011991
$ws.Range("M17").Value = 4.877755666666666
$ws.Range("N17").Value = 14.633267
$ws.Range("O17").Value = 0.09961167132870688
$ws.Range("P17").Value = 0.0996116713287069
$ws.Range("Q17").Value = 0.4529223765097777
$ws.Range("R17").Value = 4.076301388588
$ws.Range("S17").Value = 0.0004098524976074885
$ws.Range("T17").Value = 0.0004098524976074885
$ws.Range("G18").Value = 0.09285466666666665
$ws.Range("H18").Value = 0.278564
$ws.Range("I18").Value = 0.004114502770011991
$ws.Range("J18").Value = 0.004114502770011991
$ws.Range("O18").Value = 0.1360673938501395
$ws.Range("P18").Value = 0.1360673938501395
$ws.Range("Q18").Value = 0.6186821942253332
$ws.Range("R18").Value = 5.568139748028
$ws.Range("S18").Value = 0.0005598496689047114
$ws.Range("T18").Value = 0.0005598496689047115
$ws.Range("G19").Value = 0.09285466666666665
$ws.Range("H19").Value = 0.278564
$ws.Range("I19").Value = 0.004114502770011991
$ws.Range("J19").Value = 0.004114502770011991
$ws.Range("K19").Value = 3.0
$ws.Range("L19").Value = 1.0
$ws.Range("M19").Value = 0.246459
$ws.Range("N19").Value = 0.739377
$ws.Range("O19").Value = 0.00503309197542868
$ws.Range("P19").Value = 0.00503309197542868
$ws.Range("Q19").Value = 0.022884868292
$ws.Range("R19").Value = 0.205963814628
$ws.Range("S19").Value = 0.00002070867087462643
$ws.Range("T19").Value = 0.00002070867087462643
$ws.Range("G20").Value = 0.09285466666666665
$ws.Range("H20").Value = 0.278564
$ws.Range("I20").Value = 0.004114502770011991
$ws.Range("J20").Value = 0.004114502770011991
$ws.Range("M20").Value = 37.01331466666667
$ws.Range("N20").Value = 111.039944
$ws.Range("O20").Value = 0.7558718368280999
$ws.Range("P20").Value = 0.7558718368280999
$ws.Range("Q20").Value = 3.436858995601777
$ws.Range("R20").Value = 30.93173096041599
$ws.Range("S20").Value = 0.003110036766403269
$ws.Range("T20").Value = 0.003110036766403269
$ws.Range("G21").Value = 0.09285466666666665
$ws.Range("H21").Value = 0.278564
$ws.Range("I21").Value = 0.004114502770011991
$ws.Range("J21").Value = 0.004114502770011991
$ws.Range("M21").Value = 0.167274
$ws.Range("N21").Value = 0.501822
$ws.Range("O21").Value = 0.00341600601762507
$ws.Range("P21").Value = 0.00341600601762507
$ws.Range("Q21").Value = 0.015532171512
$ws.Range("R21").Value = 0.139789543608
$ws.Range("S21").Value = 0.00001405516622189598
$ws.Range("T21").Value = 0.00001405516622189598

Write-Host "Applied TPM updates to $($wb.Name)"
